$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'68.318.14"
$ws.Range("E2").Value = "'  +0.80%  "
$ws.Range("D3").Value = "'3.799.52"
$ws.Range("E3").Value = "'  -0.15%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'607.55"
$ws.Range("E5").Value = "'  +0.54%  "
$ws.Range("D6").Value = "'163.47"
$ws.Range("E6").Value = "'  -1.55%  "
$ws.Range("D7").Value = "'3.796.88"
$ws.Range("E7").Value = "'  -0.12%  "
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "'  -0.40%  "
$ws.Range("E10").Value = "'  -0.06%  "
$ws.Range("D11").Value = "'6.96"
$ws.Range("E11").Value = "'  +10.55%  "
$ws.Range("E12").Value = "'  -0.32%  "
$ws.Range("E13").Value = "'  -1.16%  "
$ws.Range("D14").Value = "'35.09"
$ws.Range("D15").Value = "'4.435.08"
$ws.Range("E15").Value = "'  -0.29%  "
$ws.Range("D16").Value = "'3.816.01"
$ws.Range("E16").Value = "'  -0.25%  "
$ws.Range("D17").Value = "'68.282.06"
$ws.Range("E17").Value = "'  +0.69%  "
$ws.Range("D18").Value = "'18.08"
$ws.Range("E18").Value = "'  -2.35%  "
$ws.Range("E19").Value = "'  +0.69%  "
$ws.Range("D20").Value = "'7.06"
$ws.Range("E20").Value = "'  -0.21%  "
$ws.Range("D21").Value = "'462.43"
$ws.Range("E21").Value = "'  -0.21%  "
$ws.Range("D22").Value = "'9.59"
$ws.Range("E22").Value = "'  -2.76%  "
$ws.Range("D23").Value = "'0.699"
$ws.Range("E23").Value = "'  -0.36%  "
$ws.Range("E24").Value = "'  +0.10%  "
$ws.Range("D25").Value = "'83.45"
$ws.Range("E25").Value = "'  +0.24%  "
$ws.Range("D26").Value = "'11.99"
$ws.Range("E26").Value = "'  -1.12%  "
$ws.Range("D27").Value = "'2.11"
$ws.Range("E27").Value = "'  -0.57%  "
$ws.Range("E28").Value = "'  -0.03%  "
$ws.Range("D29").Value = "'9.99"
$ws.Range("E29").Value = "'  -0.75%  "
$ws.Range("D30").Value = "'3.947.64"
$ws.Range("E30").Value = "'  -0.18%  "
$ws.Range("E31").Value = "'  -5.74%  "
$ws.Range("D32").Value = "'2.21"
$ws.Range("E32").Value = "'  -0.01%  "
$ws.Range("E33").Value = "'  -1.45%  "
$ws.Range("D34").Value = "'29.08"
$ws.Range("E34").Value = "'  -1.20%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "'  -0.16%  "
$ws.Range("D36").Value = "'9.03"
$ws.Range("E36").Value = "'  -0.83%  "
$ws.Range("E37").Value = "'  +0.89%  "
$ws.Range("D38").Value = "'0.149"
$ws.Range("E38").Value = "'  +7.93%  "
$ws.Range("D39").Value = "'5.87"
$ws.Range("E39").Value = "'  +1.01%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "'0.977"
$ws.Range("E40").Value = "'  -1.80%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.19"
$ws.Range("E41").Value = "'  -0.84%  "
$ws.Range("E42").Value = "'  -0.04%  "
$ws.Range("E43").Value = "'  +0.03%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'152.87"
$ws.Range("E44").Value = "'  +1.47%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.296"
$ws.Range("E45").Value = "'  -1.29%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'46.84"
$ws.Range("E46").Value = "'  -2.06%  "
$ws.Range("D47").Value = "'42.91"
$ws.Range("E47").Value = "'  -4.04%  "
$ws.Range("E48").Value = "'  +0.57%  "
$ws.Range("D49").Value = "'8.37"
$ws.Range("E49").Value = "'  +0.21%  "
$ws.Range("E50").Value = "'  +0.15%  "
$ws.Range("D51").Value = "'26.30"
$ws.Range("E51").Value = "'  -4.47%  "
